# Auto-generated Excel COM-interop edit script
# Applies the numeric cell updates described by the commit diff
# (profit/loss recompute touching columns H-N across multiple craft sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 319.9091
$ws.Range("I2").Value = 168.88889
$ws.Range("J2").Value = 999.5
$ws.Range("K2").Value = 168.88889
$ws.Range("L2").Value = 999.5
$ws.Range("M2").Value = -55.88889
$ws.Range("N2").Value = -1225.5
$ws.Range("H29").Value = 1393
$ws.Range("J29").Value = 980
$ws.Range("L29").Value = 2940
$ws.Range("N29").Value = -3502
$ws.Range("H86").Value = 1680.3636
$ws.Range("I86").Value = 1832.9445
$ws.Range("J86").Value = 993.75
$ws.Range("K86").Value = 1832.9445
$ws.Range("L86").Value = 993.75
$ws.Range("M86").Value = -709.9445000000001
$ws.Range("N86").Value = -3239.75
$ws.Range("H89").Value = 1680.3636
$ws.Range("I89").Value = 1832.9445
$ws.Range("J89").Value = 993.75
$ws.Range("K89").Value = 9164.7225
$ws.Range("L89").Value = 4968.75
$ws.Range("M89").Value = -3548.7225
$ws.Range("N89").Value = -16200.75
$ws.Range("H98").Value = 3037
$ws.Range("I98").Value = 1773.6
$ws.Range("K98").Value = 1773.6
$ws.Range("M98").Value = -275.5999999999999
$ws.Range("H122").Value = 3037
$ws.Range("I122").Value = 1773.6
$ws.Range("K122").Value = 5320.799999999999
$ws.Range("M122").Value = -2870.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 44879.5
$ws.Range("I63").Value = 50409.906
$ws.Range("J63").Value = 6166.6665
$ws.Range("K63").Value = 50409.906
$ws.Range("L63").Value = 6166.6665
$ws.Range("M63").Value = -49723.906
$ws.Range("N63").Value = -7538.6665
$ws.Range("H66").Value = 44879.5
$ws.Range("I66").Value = 50409.906
$ws.Range("J66").Value = 6166.6665
$ws.Range("K66").Value = 252049.53
$ws.Range("L66").Value = 30833.3325
$ws.Range("M66").Value = -248617.53
$ws.Range("N66").Value = -37697.3325
$ws.Range("H119").Value = 60599
$ws.Range("J119").Value = 60599
$ws.Range("L119").Value = 60599
$ws.Range("N119").Value = -70275

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1815.0834
$ws.Range("I20").Value = 846.6
$ws.Range("J20").Value = 3429.2222
$ws.Range("K20").Value = 846.6
$ws.Range("L20").Value = 3429.2222
$ws.Range("M20").Value = -599.6
$ws.Range("N20").Value = -3923.2222
$ws.Range("H86").Value = 2089.9355
$ws.Range("I86").Value = 1784.2693
$ws.Range("J86").Value = 3679.4
$ws.Range("K86").Value = 1784.2693
$ws.Range("L86").Value = 3679.4
$ws.Range("M86").Value = -661.2692999999999
$ws.Range("N86").Value = -5925.4
$ws.Range("H89").Value = 2089.9355
$ws.Range("I89").Value = 1784.2693
$ws.Range("J89").Value = 3679.4
$ws.Range("K89").Value = 8921.3465
$ws.Range("L89").Value = 18397
$ws.Range("M89").Value = -3305.3465
$ws.Range("N89").Value = -29629
$ws.Range("H107").Value = 3545.1904
$ws.Range("I107").Value = 1743.6154
$ws.Range("J107").Value = 6472.75
$ws.Range("K107").Value = 1743.6154
$ws.Range("L107").Value = 6472.75
$ws.Range("M107").Value = 176.3846000000001
$ws.Range("N107").Value = -10312.75
$ws.Range("H134").Value = 9800.629999999999
$ws.Range("I134").Value = 6436.2
$ws.Range("K134").Value = 19308.6
$ws.Range("M134").Value = -16773.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2673.2354
$ws.Range("I62").Value = 2183.1538
$ws.Range("K62").Value = 2183.1538
$ws.Range("M62").Value = -1559.1538
$ws.Range("H65").Value = 2673.2354
$ws.Range("I65").Value = 2183.1538
$ws.Range("K65").Value = 10915.769
$ws.Range("M65").Value = -7795.769
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H141").Value = 337499
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 2151.682
$ws.Range("J38").Value = 3887
$ws.Range("L38").Value = 11661
$ws.Range("N38").Value = -12355
$ws.Range("H86").Value = 839.1875
$ws.Range("J86").Value = 830.8333
$ws.Range("L86").Value = 2492.4999
$ws.Range("N86").Value = -4864.4999
$ws.Range("H87").Value = 8432
$ws.Range("I87").Value = 8432
$ws.Range("K87").Value = 25296
$ws.Range("M87").Value = -24048
$ws.Range("H89").Value = 839.1875
$ws.Range("J89").Value = 830.8333
$ws.Range("L89").Value = 7477.4997
$ws.Range("N89").Value = -19333.4997
$ws.Range("H90").Value = 8432
$ws.Range("I90").Value = 8432
$ws.Range("K90").Value = 75888
$ws.Range("M90").Value = -69648

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1957.4
$ws.Range("J3").Value = 1749.5
$ws.Range("L3").Value = 1749.5
$ws.Range("N3").Value = -1981.5
$ws.Range("H10").Value = 2691.625
$ws.Range("I10").Value = 2790.4285
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 2790.4285
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = -2621.4285
$ws.Range("N10").Value = -2338
$ws.Range("H31").Value = 1417
$ws.Range("I31").Value = 1417
$ws.Range("K31").Value = 1417
$ws.Range("M31").Value = -1125
$ws.Range("H37").Value = 1417
$ws.Range("I37").Value = 1417
$ws.Range("K37").Value = 1417
$ws.Range("M37").Value = -1140
$ws.Range("H40").Value = 29666.334
$ws.Range("I40").Value = 29499.5
$ws.Range("K40").Value = 29499.5
$ws.Range("M40").Value = -29348.5
$ws.Range("H70").Value = 8755.4
$ws.Range("I70").Value = 7932.6665
$ws.Range("J70").Value = 9989.5
$ws.Range("K70").Value = 7932.6665
$ws.Range("L70").Value = 9989.5
$ws.Range("M70").Value = -7662.6665
$ws.Range("N70").Value = -10529.5
$ws.Range("H73").Value = 8755.4
$ws.Range("I73").Value = 7932.6665
$ws.Range("J73").Value = 9989.5
$ws.Range("K73").Value = 7932.6665
$ws.Range("L73").Value = 9989.5
$ws.Range("M73").Value = -6996.6665
$ws.Range("N73").Value = -11861.5
$ws.Range("H80").Value = 10553.35
$ws.Range("I80").Value = 11010.625
$ws.Range("J80").Value = 10248.5
$ws.Range("K80").Value = 11010.625
$ws.Range("L80").Value = 10248.5
$ws.Range("M80").Value = -10012.625
$ws.Range("N80").Value = -12244.5
$ws.Range("H83").Value = 10553.35
$ws.Range("I83").Value = 11010.625
$ws.Range("J83").Value = 10248.5
$ws.Range("K83").Value = 55053.125
$ws.Range("L83").Value = 51242.5
$ws.Range("M83").Value = -50061.125
$ws.Range("N83").Value = -61226.5
$ws.Range("H92").Value = 34996.8
$ws.Range("J92").Value = 34996.8
$ws.Range("L92").Value = 34996.8
$ws.Range("N92").Value = -38740.8
$ws.Range("H102").Value = 1709.9231
$ws.Range("I102").Value = 1575.5555
$ws.Range("K102").Value = 1575.5555
$ws.Range("M102").Value = 46.44450000000006
$ws.Range("H113").Value = 2248.8
$ws.Range("I113").Value = 2248.8
$ws.Range("K113").Value = 2248.8
$ws.Range("M113").Value = -78.80000000000018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1122.4736
$ws.Range("I46").Value = 1113.6666
$ws.Range("K46").Value = 1113.6666
$ws.Range("M46").Value = -925.6666
$ws.Range("H119").Value = 80672
$ws.Range("J119").Value = 80672
$ws.Range("L119").Value = 80672
$ws.Range("N119").Value = -90348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2996
$ws.Range("I3").Value = 2992.5
$ws.Range("J3").Value = 2999.5
$ws.Range("K3").Value = 2992.5
$ws.Range("L3").Value = 2999.5
$ws.Range("M3").Value = -2878.5
$ws.Range("N3").Value = -3227.5
$ws.Range("H10").Value = 49999.332
$ws.Range("I10").Value = 49999.332
$ws.Range("K10").Value = 49999.332
$ws.Range("M10").Value = -49830.332
$ws.Range("H62").Value = 2866.5
$ws.Range("J62").Value = 2859.8
$ws.Range("L62").Value = 2859.8
$ws.Range("N62").Value = -4107.8
$ws.Range("H65").Value = 2866.5
$ws.Range("J65").Value = 2859.8
$ws.Range("L65").Value = 14299
$ws.Range("N65").Value = -20539
$ws.Range("H119").Value = 90924.25
$ws.Range("J119").Value = 90924.25
$ws.Range("L119").Value = 90924.25
$ws.Range("N119").Value = -100600.25
$ws.Range("H136").Value = 3411.879
$ws.Range("I136").Value = 3203.6428
$ws.Range("K136").Value = 9610.928400000001
$ws.Range("M136").Value = -7060.928400000001
